$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 110; existing rows 110-120 shift down to 111-121.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new weekly price-report entry
# (same Mercado / Región / Producto metadata as the surrounding rows).
$ws.Range("A110").Value = 10
$ws.Range("B110").Value = "Vega Modelo de Temuco"
$ws.Range("C110").Value = "La Araucanía"
$ws.Range("D110").Value = 45180
$ws.Range("E110").Value = 9
$ws.Range("F110").Value = "Fruta"
$ws.Range("G110").Value = 100108
$ws.Range("H110").Value = "Tropicales y subtropicales"
$ws.Range("I110").Value = 100108003
$ws.Range("J110").Value = "Maracuyá"
$ws.Range("K110").Value = "Sin especificar"
$ws.Range("L110").Value = "Primera"
$ws.Range("M110").Value = 80
$ws.Range("N110").Value = 50000
$ws.Range("O110").Value = 50000
$ws.Range("P110").Value = 50000
$ws.Range("Q110").Value = "$/caja 18 kilos"
$ws.Range("R110").Value = "Región de Arica y Parinacota"
$ws.Range("S110").Value = 2778
$ws.Range("T110").Value = 18
